$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 279.45947
$ws.Range("I33").Value = 300.72726
$ws.Range("J33").Value = 104
$ws.Range("K33").Value = 300.72726
$ws.Range("L33").Value = 104
$ws.Range("M33").Value = -71.72726
$ws.Range("N33").Value = -562
$ws.Range("H40").Value = 1758.3334
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 1781.8182
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1781.8182
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -2131.8182
$ws.Range("H64").Value = 3360.7693
$ws.Range("I64").Value = 3100
$ws.Range("J64").Value = 3523.75
$ws.Range("K64").Value = 3100
$ws.Range("L64").Value = 3523.75
$ws.Range("M64").Value = -2852
$ws.Range("N64").Value = -4019.75
$ws.Range("H67").Value = 3360.7693
$ws.Range("I67").Value = 3100
$ws.Range("J67").Value = 3523.75
$ws.Range("K67").Value = 3100
$ws.Range("L67").Value = 3523.75
$ws.Range("M67").Value = -2242
$ws.Range("N67").Value = -5239.75
$ws.Range("H116").Value = 1990.5
$ws.Range("I116").Value = 1983.1666
$ws.Range("J116").Value = 2001.5
$ws.Range("K116").Value = 1983.1666
$ws.Range("L116").Value = 2001.5
$ws.Range("M116").Value = 1458.8334
$ws.Range("N116").Value = -8885.5
$ws.Range("H132").Value = 5496590.5
$ws.Range("I132").Value = 6495606
$ws.Range("J132").Value = 2004.75
$ws.Range("K132").Value = 19486818
$ws.Range("L132").Value = 6014.25
$ws.Range("M132").Value = -19484288
$ws.Range("N132").Value = -11074.25
$ws.Range("H138").Value = 1958.8572
$ws.Range("I138").Value = 1515.3793
$ws.Range("J138").Value = 2435.1853
$ws.Range("K138").Value = 4546.1379
$ws.Range("L138").Value = 7305.5559
$ws.Range("M138").Value = 593.8621000000003
$ws.Range("N138").Value = -17585.5559

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4995.8335
$ws.Range("I32").Value = 4070.712
$ws.Range("J32").Value = 9194.462
$ws.Range("K32").Value = 4070.712
$ws.Range("L32").Value = 9194.462
$ws.Range("M32").Value = -3783.712
$ws.Range("N32").Value = -9768.462
$ws.Range("H88").Value = 1116702.2
$ws.Range("I88").Value = 2508251.5
$ws.Range("K88").Value = 2508251.5
$ws.Range("M88").Value = -2507845.5
$ws.Range("H91").Value = 1116702.2
$ws.Range("I91").Value = 2508251.5
$ws.Range("K91").Value = 2508251.5
$ws.Range("M91").Value = -2506847.5
$ws.Range("H115").Value = 44742
$ws.Range("J115").Value = 44742
$ws.Range("L115").Value = 44742
$ws.Range("N115").Value = -47876
$ws.Range("H132").Value = 5751.6772
$ws.Range("I132").Value = 5972.654
$ws.Range("J132").Value = 4602.6
$ws.Range("K132").Value = 17917.962
$ws.Range("L132").Value = 13807.8
$ws.Range("M132").Value = -15387.962
$ws.Range("N132").Value = -18867.8

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4087.4
$ws.Range("I86").Value = 2868
$ws.Range("J86").Value = 5916.5
$ws.Range("K86").Value = 2868
$ws.Range("L86").Value = 5916.5
$ws.Range("M86").Value = -1745
$ws.Range("N86").Value = -8162.5
$ws.Range("H89").Value = 4087.4
$ws.Range("I89").Value = 2868
$ws.Range("J89").Value = 5916.5
$ws.Range("K89").Value = 14340
$ws.Range("L89").Value = 29582.5
$ws.Range("M89").Value = -8724
$ws.Range("N89").Value = -40814.5
$ws.Range("H134").Value = 32109.766
$ws.Range("I134").Value = 55179.156
$ws.Range("J134").Value = 2888.5334
$ws.Range("K134").Value = 165537.468
$ws.Range("L134").Value = 8665.600199999999
$ws.Range("M134").Value = -163002.468
$ws.Range("N134").Value = -13735.6002

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2527204.8
$ws.Range("I31").Value = 1474.7021
$ws.Range("J31").Value = 8775063
$ws.Range("K31").Value = 1474.7021
$ws.Range("L31").Value = 8775063
$ws.Range("M31").Value = -1179.7021
$ws.Range("N31").Value = -8775653
$ws.Range("H34").Value = 2527204.8
$ws.Range("I34").Value = 1474.7021
$ws.Range("J34").Value = 8775063
$ws.Range("K34").Value = 1474.7021
$ws.Range("L34").Value = 8775063
$ws.Range("M34").Value = -1272.7021
$ws.Range("N34").Value = -8775467
$ws.Range("H62").Value = 111113016
$ws.Range("J62").Value = 111113016
$ws.Range("L62").Value = 111113016
$ws.Range("N62").Value = -111114264
$ws.Range("H65").Value = 111113016
$ws.Range("J65").Value = 111113016
$ws.Range("L65").Value = 555565080
$ws.Range("N65").Value = -555571320
$ws.Range("H123").Value = 30390
$ws.Range("J123").Value = 30390
$ws.Range("L123").Value = 30390
$ws.Range("N123").Value = -40190
$ws.Range("H132").Value = 4086.9333
$ws.Range("I132").Value = 3288.5
$ws.Range("J132").Value = 4999.4287
$ws.Range("K132").Value = 9865.5
$ws.Range("L132").Value = 14998.2861
$ws.Range("M132").Value = -7335.5
$ws.Range("N132").Value = -20058.2861
$ws.Range("H134").Value = 1060.7142
$ws.Range("I134").Value = 945.7727
$ws.Range("J134").Value = 1482.1666
$ws.Range("K134").Value = 2837.3181
$ws.Range("L134").Value = 4446.4998
$ws.Range("M134").Value = -302.3181
$ws.Range("N134").Value = -9516.4998

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 33293.312
$ws.Range("I132").Value = 36047.863
$ws.Range("J132").Value = 6666
$ws.Range("K132").Value = 108143.589
$ws.Range("L132").Value = 19998
$ws.Range("M132").Value = -105613.589
$ws.Range("N132").Value = -25058

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 31251582
$ws.Range("I7").Value = 1211.4546
$ws.Range("J7").Value = 100002400
$ws.Range("K7").Value = 1211.4546
$ws.Range("L7").Value = 100002400
$ws.Range("M7").Value = -1099.4546
$ws.Range("N7").Value = -100002624
$ws.Range("H40").Value = 2217.375
$ws.Range("I40").Value = 1823.1666
$ws.Range("K40").Value = 1823.1666
$ws.Range("M40").Value = -1687.1666
$ws.Range("H46").Value = 1238.0435
$ws.Range("I46").Value = 296.16666
$ws.Range("J46").Value = 1570.4706
$ws.Range("K46").Value = 296.16666
$ws.Range("L46").Value = 1570.4706
$ws.Range("M46").Value = -108.16666
$ws.Range("N46").Value = -1946.4706
$ws.Range("H122").Value = 2800
$ws.Range("J122").Value = 2888.889
$ws.Range("L122").Value = 8666.667000000001
$ws.Range("N122").Value = -13566.667
$ws.Range("H126").Value = 31251582
$ws.Range("I126").Value = 1211.4546
$ws.Range("J126").Value = 100002400
$ws.Range("K126").Value = 3634.3638
$ws.Range("L126").Value = 300007200
$ws.Range("M126").Value = -1164.3638
$ws.Range("N126").Value = -300012140
$ws.Range("H132").Value = 2519.682
$ws.Range("I132").Value = 2267.5
$ws.Range("J132").Value = 2961
$ws.Range("K132").Value = 6802.5
$ws.Range("L132").Value = 8883
$ws.Range("M132").Value = -4272.5
$ws.Range("N132").Value = -13943

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3287.5
$ws.Range("I14").Value = 3322
$ws.Range("J14").Value = 3230
$ws.Range("K14").Value = 3322
$ws.Range("L14").Value = 3230
$ws.Range("M14").Value = -3154
$ws.Range("N14").Value = -3566
$ws.Range("H81").Value = 1842.6666
$ws.Range("I81").Value = 1611.2
$ws.Range("K81").Value = 3222.4
$ws.Range("M81").Value = -2161.4
$ws.Range("H84").Value = 1842.6666
$ws.Range("I84").Value = 1611.2
$ws.Range("K84").Value = 16112
$ws.Range("M84").Value = -10808
$ws.Range("H113").Value = 573.5
$ws.Range("J113").Value = 416.5
$ws.Range("L113").Value = 1249.5
$ws.Range("N113").Value = -5589.5
$ws.Range("H132").Value = 1315.4531
$ws.Range("I132").Value = 1323.6735
$ws.Range("J132").Value = 1288.6
$ws.Range("K132").Value = 3971.020500000001
$ws.Range("L132").Value = 3865.8
$ws.Range("M132").Value = -1441.020500000001
$ws.Range("N132").Value = -8925.799999999999
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()
